# Add a new worksheet "Treasure For CR 5 Coin" at the end of the workbook,
# seeded from the existing "Treasure For CR 0 Coin" sheet (same layout/styles),
# but with a "bad" die designation ("3D20" instead of "D20") in B1 so the new
# internal function _check_roll_column() has a bad-roll-column case to catch.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Treasure For CR 0 Coin")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$new = $wb.Worksheets.Add($null, $lastSheet)
$new.Name = "Treasure For CR 5 Coin"

# Copy layout, values, and formatting from the CR 0 Coin sheet wholesale ...
$src.Range("A1:C6").Copy($new.Range("A1"))

# ... then poison the die designation to produce the "bad roll column" test case.
$new.Range("B1").Value = "3D20"

# The copy leaves A1 blank (it was blank on the source sheet too) plus a
# couple of other blank-but-touched cells; clear them back out so the sheet
# matches the source sheet's actually-used range exactly.
$new.Range("A1").Clear()
$new.Range("C2").Clear()
$new.Range("B6:C6").Clear()

# Make the new sheet the active one, with the same lingering selection the
# author left behind.
$new.Activate()
$new.Range("C21").Select()
